$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVTQaZ")

# Remove plugin hybrid vehicles (row 6) as qualifying ZEVs for 2020-2029 (columns B:K)
$ws.Range("B6:K6").Value = 0

# Update selection to reflect where the edit left the cursor
$ws.Activate()
$ws.Range("D14").Select()
